$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C8").Value = 123209
$ws.Range("E8").Value = 267472682
$ws.Range("C10").Value = 139393
$ws.Range("D10").Value = 28151
$ws.Range("E10").Value = 530548909
$ws.Range("C13").Value = 23698
$ws.Range("D13").Value = 7768
$ws.Range("E13").Value = 44025148
$ws.Range("C14").Value = 75728
$ws.Range("E14").Value = 149750461
$ws.Range("C18").Value = 54256
$ws.Range("D18").Value = 23500
$ws.Range("E18").Value = 72581991
$ws.Range("C33").Value = 14319
$ws.Range("E33").Value = 28709605
$ws.Range("C52").Value = 21512
$ws.Range("E52").Value = 45371495
$ws.Range("C54").Value = 11840
$ws.Range("D54").Value = 3062
$ws.Range("E54").Value = 19139113
$ws.Range("C55").Value = 19376
$ws.Range("E55").Value = 25932741
$ws.Range("C63").Value = 29399
$ws.Range("D63").Value = 9593
$ws.Range("E63").Value = 59320718
$ws.Range("C69").Value = 13449
$ws.Range("E69").Value = 26547172
$ws.Range("C156").Value = 30407
$ws.Range("D156").Value = 9557
$ws.Range("E156").Value = 58879753
$ws.Range("C162").Value = 47589
$ws.Range("E162").Value = 78214579
$ws.Range("C169").Value = 360938
$ws.Range("E169").Value = 723181594
$ws.Range("D174").Value = 65757
$ws.Range("E174").Value = 418616368
$ws.Range("C179").Value = 137239
$ws.Range("D179").Value = 27285
$ws.Range("E179").Value = 280117395
$ws.Range("C239").Value = 42834
$ws.Range("E239").Value = 156268843
$ws.Range("C243").Value = 17900
$ws.Range("D243").Value = 5673
$ws.Range("E243").Value = 37254909
$ws.Range("C248").Value = 29598
$ws.Range("E248").Value = 51916872
$ws.Range("C261").Value = 50834
$ws.Range("D261").Value = 16405
$ws.Range("E261").Value = 99325285
$ws.Range("C265").Value = 38786
$ws.Range("D265").Value = 17017
$ws.Range("E265").Value = 51106449
$ws.Range("C266").Value = 40450
$ws.Range("D266").Value = 8474
$ws.Range("E266").Value = 78774364
$ws.Range("C267").Value = 66601
$ws.Range("E267").Value = 111884409
$ws.Range("C274").Value = 114390
$ws.Range("D274").Value = 34563
$ws.Range("E274").Value = 227403209
$ws.Range("C276").Value = 112794
$ws.Range("E276").Value = 389927023
$ws.Range("C283").Value = 38106
$ws.Range("D283").Value = 9584
$ws.Range("E283").Value = 54837800
$ws.Range("C284").Value = 43153
$ws.Range("D284").Value = 19190
$ws.Range("E284").Value = 57104588
$ws.Range("C285").Value = 48380
$ws.Range("D285").Value = 9714
$ws.Range("E285").Value = 95026612
$ws.Range("C299").Value = 27336
$ws.Range("E299").Value = 56885524
$ws.Range("C300").Value = 9991
$ws.Range("E300").Value = 21694176
$ws.Range("C304").Value = 34291
$ws.Range("D304").Value = 10133
$ws.Range("E304").Value = 62643914
$ws.Range("C317").Value = 65759
$ws.Range("E317").Value = 127218841
$ws.Range("C318").Value = 37003
$ws.Range("E318").Value = 76892666
$ws.Range("C322").Value = 46048
$ws.Range("D322").Value = 9120
$ws.Range("E322").Value = 94697981
